# Generate Report for Handoff
# Updates the "Priority" column to "ht" for the rows that are ready for
# handoff, and bumps the related timestamp columns forward by 15 seconds
# to reflect the newly generated handoff report.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 12, 13)

# zh-cn sheet: Priority column (E) -> "ht", Latest Handoff Datetime (H) +15s
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-30 06:23:43"
}

# de-de sheet: Priority column (E) -> "ht", Latest Handoff Datetime (H) +15s
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-30 06:23:48"
}
# (Overview.G and de-de.H share the same new timestamp: 2016-08-30 06:23:48)

# Overview sheet: Latest HO Xliff Generate Date (G) +15s
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-30 06:23:48"
}
